$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.800.87'
$ws.Range("E2").Value = '  +1.74%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.477.05'
$ws.Range("E3").Value = '  +1.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.51'
$ws.Range("E5").Value = '  +1.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.11'
$ws.Range("E6").Value = '  +2.54%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +1.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.474.87'
$ws.Range("E9").Value = '  +1.49%  '

$ws.Range("E10").Value = '  +1.28%  '

$ws.Range("E11").Value = '  +1.04%  '

$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.360'
$ws.Range("E12").Value = '  +1.78%  '

$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.28'
$ws.Range("E13").Value = '  +0.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.35'
$ws.Range("E14").Value = '  +2.02%  '

$ws.Range("E15").Value = '  -1.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.926.37'
$ws.Range("E16").Value = '  +1.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.786.98'
$ws.Range("E17").Value = '  +2.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.478.67'
$ws.Range("E18").Value = '  +1.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.57'
$ws.Range("E19").Value = '  +2.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.42'
$ws.Range("E20").Value = '  +6.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.23'
$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.13'
$ws.Range("E23").Value = '  +19.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.05'
$ws.Range("E25").Value = '  -1.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '629.96'
$ws.Range("E26").Value = '  +11.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000106'
$ws.Range("E27").Value = '  +3.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.79'
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.54'
$ws.Range("E29").Value = '  +5.74%  '

$ws.Range("B30").Value = 'WrappedeETH'
$ws.Range("C30").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.599.36'
$ws.Range("E30").Value = '  +1.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.39'
$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("E33").Value = '  -2.34%  '

$ws.Range("E34").Value = '  +2.56%  '

$ws.Range("E35").Value = '  +8.56%  '

$ws.Range("E36").Value = '  -0.41%  '

$ws.Range("E37").Value = '  +0.09%  '

$ws.Range("E38").Value = '  +0.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.51'
$ws.Range("E39").Value = '  +1.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.89'
$ws.Range("E40").Value = '  +0.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.75'
$ws.Range("E41").Value = '  +13.20%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.83'
$ws.Range("E42").Value = '  -0.50%  '

$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '147.47'
$ws.Range("E43").Value = '  -0.64%  '

$ws.Range("E44").Value = '  -0.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '151.27'
$ws.Range("E45").Value = '  +1.84%  '

$ws.Range("E46").Value = '  +3.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.51'
$ws.Range("E47").Value = '  +4.65%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0545'
$ws.Range("E48").Value = '  +1.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.607'
$ws.Range("E49").Value = '  +0.98%  '

$ws.Range("E50").Value = '  +2.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0923'
$ws.Range("E51").Value = '  -0.54%  '
